$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D23").Value = "cusolver64_10.dll not found만 찾을 수 없다는 에러"
$ws.Range("E23").Value = "https://theonly1.tistory.com/2756"

$ws.Range("D39").Value = "How to Think Probabilistically with Discrete Distributions"
$ws.Range("E39").Value = "https://a292run.tistory.com/entry/How-to-Think-Probabilistically-with-Discrete-Distributions-1"

$ws.Range("D43").Value = "iptime 남는 공유기 증폭기 확장기로 쓰기"
$ws.Range("E43").Value = "https://nittaku.tistory.com/507"

$ws.Range("D46").Value = "[한국생명공학연구원] 2021년 03월, 생물정보학(Bioinformatics 채용), 생명정보분석프로그램개발, 바이오데이터품질관리, 인체유래데이터수집"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/388"

$ws.Range("D51").Value = "[github] 자주 사용하는 마크다운(markdown) 문법 정리"
$ws.Range("E51").Value = "https://bskyvision.com/1140"
